$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new row 50 with raw/clean data for 2020-07-19.
# Column A holds the date as text (matching the existing rows, which are
# shared strings rather than real dates), so we build it via a text
# formula and convert it to a static value/shared-string in place -
# this avoids Excel's automatic text->date parsing while also avoiding
# any style/number-format side effects.
$ws.Cells.Item(50, 1).Formula = '="2020-07-19"'
$ws.Cells.Item(50, 1).Copy()
$ws.Cells.Item(50, 1).PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$ws.Cells.Item(50, 2).Value = 344224
$ws.Cells.Item(50, 3).Value = 394156
$ws.Cells.Item(50, 4).Value = 83542
$ws.Cells.Item(50, 5).Value = 39184
$ws.Cells.Item(50, 6).Value = 28.55
